# Reduce test file size: the sheet previously held the full exported
# dataset (300k+ rows via the sortState range) which required Git LFS.
# Replace it with a single representative data row so the test fixture
# stays small while still exercising the expected shape (year 2000).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2000

# Matches the author's resulting selection/cursor position after the edit.
$ws.Range("A3").Select()

$wb.Save()
